$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.73
$ws.Range("R9").Value = 2.1
$ws.Range("U9").Value = 2.75
$ws.Range("V9").Value = 1.44
$ws.Range("Y9").Value = 1.62
$ws.Range("Z9").Value = 2.2
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 7
$ws.Range("AK9").Value = 151
$ws.Range("AL9").Value = 12

# Row 15
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 1.67
$ws.Range("AG15").Value = 8.5
$ws.Range("AN15").Value = 9.5

# Row 20
$ws.Range("Q20").Value = 2.1
$ws.Range("R20").Value = 1.7

# Row 24
$ws.Range("G24").Value = 3
$ws.Range("I24").Value = 2.35
$ws.Range("AA24").Value = 8
$ws.Range("AC24").Value = 11
$ws.Range("AE24").Value = 26
$ws.Range("AM24").Value = 11
$ws.Range("AO24").Value = 23

# Row 25
$ws.Range("G25").Value = 4.75
$ws.Range("H25").Value = 3.3
$ws.Range("I25").Value = 1.85
$ws.Range("K25").Value = 1.91
$ws.Range("L25").Value = 2.63
$ws.Range("M25").Value = 1.13
$ws.Range("N25").Value = 6
$ws.Range("O25").Value = 1.62
$ws.Range("P25").Value = 2.2
$ws.Range("Q25").Value = 2.88
$ws.Range("R25").Value = 1.4
$ws.Range("S25").Value = 5
$ws.Range("T25").Value = 1.18
$ws.Range("U25").Value = 6
$ws.Range("V25").Value = 1.13
$ws.Range("W25").Value = 1.62
$ws.Range("X25").Value = 2.2
$ws.Range("Y25").Value = 2.5
$ws.Range("Z25").Value = 1.5
$ws.Range("AA25").Value = 8.5
$ws.Range("AB25").Value = 21
$ws.Range("AE25").Value = 51
$ws.Range("AG25").Value = 6
$ws.Range("AH25").Value = 6.5
$ws.Range("AJ25").Value = 101
$ws.Range("AL25").Value = 4.75
$ws.Range("AO25").Value = 15
$ws.Range("AP25").Value = 21
$ws.Range("AR25").Value = 2.1
$ws.Range("AS25").Value = 1.78

# Row 26
$ws.Range("G26").Value = 1.95
$ws.Range("I26").Value = 4.33
$ws.Range("J26").Value = 2.75
$ws.Range("L26").Value = 5.5
$ws.Range("M26").Value = 1.11
$ws.Range("N26").Value = 6.5
$ws.Range("Q26").Value = 2.7
$ws.Range("R26").Value = 1.44
$ws.Range("S26").Value = 4.8
$ws.Range("T26").Value = 1.19
$ws.Range("AB26").Value = 7.5
$ws.Range("AD26").Value = 17
$ws.Range("AH26").Value = 6.5
$ws.Range("AI26").Value = 23
$ws.Range("AN26").Value = 17
$ws.Range("AO26").Value = 51

# Row 28
$ws.Range("M28").Value = 1.1
$ws.Range("N28").Value = 7
$ws.Range("Q28").Value = 2.6
$ws.Range("R28").Value = 1.48
$ws.Range("U28").Value = 5.5
$ws.Range("V28").Value = 1.14
$ws.Range("AR28").Value = 1.98
$ws.Range("AS28").Value = 1.88

# Row 29
$ws.Range("AR29").Value = 2
$ws.Range("AS29").Value = 1.85

# Row 39
$ws.Range("G39").Value = 1.45
$ws.Range("H39").Value = 4.1
$ws.Range("I39").Value = 7
$ws.Range("L39").Value = 7.5
$ws.Range("M39").Value = 1.06
$ws.Range("N39").Value = 10
$ws.Range("Q39").Value = 2.15
$ws.Range("R39").Value = 1.67
$ws.Range("S39").Value = 3.25
$ws.Range("T39").Value = 1.33
$ws.Range("Y39").Value = 2.5
$ws.Range("Z39").Value = 1.5
$ws.Range("AA39").Value = 5
$ws.Range("AD39").Value = 9.5
$ws.Range("AF39").Value = 41
$ws.Range("AH39").Value = 8.5
$ws.Range("AI39").Value = 26
$ws.Range("AJ39").Value = 101
$ws.Range("AL39").Value = 13
$ws.Range("AM39").Value = 34
$ws.Range("AQ39").Value = 67
$ws.Range("AR39").Value = 1.64
$ws.Range("AS39").Value = 2.22

# Row 40
$ws.Range("J40").Value = 1.8

# Row 41
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 2.88
$ws.Range("I41").Value = 2.6
$ws.Range("J41").Value = 3.75
$ws.Range("L41").Value = 3.4
$ws.Range("AD41").Value = 29
$ws.Range("AJ41").Value = 67
$ws.Range("AL41").Value = 6.5
$ws.Range("AM41").Value = 11
$ws.Range("AN41").Value = 11
$ws.Range("AP41").Value = 26

# Row 77
$ws.Range("G77").Value = 2.75
$ws.Range("H77").Value = 3.3
$ws.Range("L77").Value = 3.25
$ws.Range("O77").Value = 1.33
$ws.Range("P77").Value = 3.25
$ws.Range("Q77").Value = 2.08
$ws.Range("R77").Value = 1.73
$ws.Range("U77").Value = 3.75
$ws.Range("V77").Value = 1.25
$ws.Range("W77").Value = 1.44
$ws.Range("X77").Value = 2.63
$ws.Range("Y77").Value = 1.8
$ws.Range("Z77").Value = 1.95
$ws.Range("AA77").Value = 8.5
$ws.Range("AE77").Value = 23
$ws.Range("AF77").Value = 34
$ws.Range("AG77").Value = 9.5
$ws.Range("AK77").Value = 251
$ws.Range("AL77").Value = 8

# Row 86
$ws.Range("H86").Value = 7
$ws.Range("I86").Value = 1.2
$ws.Range("J86").Value = 11
$ws.Range("O86").Value = 1.14
$ws.Range("P86").Value = 5.5
$ws.Range("Q86").Value = 1.5
$ws.Range("R86").Value = 2.5
$ws.Range("S86").Value = 1.83
$ws.Range("T86").Value = 2.03
$ws.Range("W86").Value = 1.25
$ws.Range("X86").Value = 3.75
$ws.Range("Y86").Value = 2.2
$ws.Range("Z86").Value = 1.62
$ws.Range("AA86").Value = 26
$ws.Range("AC86").Value = 34
$ws.Range("AE86").Value = 81

# Row 102
$ws.Range("G102").Value = 2.1
$ws.Range("H102").Value = 3.1
$ws.Range("I102").Value = 3.8
$ws.Range("J102").Value = 2.88
$ws.Range("K102").Value = 1.95
$ws.Range("L102").Value = 4.5
$ws.Range("W102").Value = 1.53
$ws.Range("X102").Value = 2.38
$ws.Range("AC102").Value = 9.5
$ws.Range("AD102").Value = 19
$ws.Range("AL102").Value = 9
$ws.Range("AN102").Value = 15

# Row 119
$ws.Range("G119").Value = 2.4
$ws.Range("H119").Value = 2.9
$ws.Range("I119").Value = 3.05
$ws.Range("J119").Value = 2.9
$ws.Range("K119").Value = 2.02
$ws.Range("L119").Value = 3.55
$ws.Range("P119").Value = 2.77
$ws.Range("W119").Value = 1.4
$ws.Range("X119").Value = 2.52
$ws.Range("Y119").Value = 1.7
$ws.Range("Z119").Value = 1.93
$ws.Range("AA119").Value = 7.7
$ws.Range("AB119").Value = 12
$ws.Range("AC119").Value = 9
$ws.Range("AD119").Value = 26
$ws.Range("AF119").Value = 29
$ws.Range("AG119").Value = 8.25
$ws.Range("AH119").Value = 5.7
$ws.Range("AI119").Value = 13
$ws.Range("AJ119").Value = 60
$ws.Range("AL119").Value = 8.75
$ws.Range("AN119").Value = 10.75
$ws.Range("AP119").Value = 28
$ws.Range("AQ119").Value = 35

# Row 120
$ws.Range("G120").Value = 2.47
$ws.Range("H120").Value = 2.55
$ws.Range("I120").Value = 3.4
$ws.Range("J120").Value = 3.15
$ws.Range("K120").Value = 1.82
$ws.Range("L120").Value = 4
$ws.Range("M120").Value = 1.1
$ws.Range("N120").Value = 6.4
$ws.Range("O120").Value = 1.47
$ws.Range("P120").Value = 2.32
$ws.Range("Q120").Value = 2.37
$ws.Range("R120").Value = 1.45
$ws.Range("W120").Value = 1.53
$ws.Range("X120").Value = 2.2
$ws.Range("Y120").Value = 1.87
$ws.Range("Z120").Value = 1.75
$ws.Range("AA120").Value = 6.3
$ws.Range("AB120").Value = 11.5
$ws.Range("AC120").Value = 9.5
$ws.Range("AD120").Value = 29
$ws.Range("AE120").Value = 24
$ws.Range("AG120").Value = 5.9
$ws.Range("AH120").Value = 5.1
$ws.Range("AI120").Value = 14.5
$ws.Range("AJ120").Value = 80
$ws.Range("AK120").Value = 800
$ws.Range("AL120").Value = 7.9
$ws.Range("AM120").Value = 17.5
$ws.Range("AN120").Value = 11.75
$ws.Range("AO120").Value = 55
$ws.Range("AQ120").Value = 45

# Row 121
$ws.Range("G121").Value = 1.72
$ws.Range("H121").Value = 3.15
$ws.Range("I121").Value = 4.6
$ws.Range("J121").Value = 2.32
$ws.Range("K121").Value = 2.05
$ws.Range("O121").Value = 1.42
$ws.Range("P121").Value = 2.75
$ws.Range("Q121").Value = 2.29
$ws.Range("R121").Value = 1.58
$ws.Range("U121").Value = 4.2
$ws.Range("V121").Value = 1.2
$ws.Range("W121").Value = 1.49
$ws.Range("X121").Value = 2.45
$ws.Range("Y121").Value = 2.12
$ws.Range("Z121").Value = 1.64
$ws.Range("AA121").Value = 4.4
$ws.Range("AD121").Value = 10
$ws.Range("AE121").Value = 13
$ws.Range("AG121").Value = 7.4
$ws.Range("AH121").Value = 5.8
$ws.Range("AI121").Value = 16
$ws.Range("AL121").Value = 8.4
$ws.Range("AM121").Value = 21
$ws.Range("AN121").Value = 14
